# Generate Report for Handoff
#
# Updates the "3220368b-681a-4b12-811d-c714053ebfb2.md" row (row 4) across the
# Overview, zh-cn and de-de sheets to reflect a freshly generated handoff
# report: the per-language "Latest Handoff Datetime" columns get a new,
# later timestamp, and the Overview sheet's "Latest HO Xliff Generate Date"
# is refreshed to the newest of those timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn: Latest Handoff Datetime (column H) for 3220368b-681a-4b12-811d-c714053ebfb2.md (row 4)
$wsZhCn.Range("H4").Value = "2016-10-27 08:51:18"

# de-de: Latest Handoff Datetime (column H) for 3220368b-681a-4b12-811d-c714053ebfb2.md (row 4)
$wsDeDe.Range("H4").Value = "2016-10-27 08:51:29"

# Overview: Latest HO Xliff Generate Date (column G) for 3220368b-681a-4b12-811d-c714053ebfb2.md (row 4)
$wsOverview.Range("G4").Value = "2016-10-27 08:51:29"
